$wb = $excel.ActiveWorkbook

# This script applies a market-data refresh to columns H-N (price/profit columns)
# across multiple worksheets, per the scheduled runner update.

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 33667
$ws.Range("I12").Value = 33667
$ws.Range("K12").Value = 33667
$ws.Range("M12").Value = -33497
$ws.Range("H41").Value = 1796.3334
$ws.Range("I41").Value = 366.33334
$ws.Range("J41").Value = 2273
$ws.Range("K41").Value = 366.33334
$ws.Range("L41").Value = 2273
$ws.Range("M41").Value = 73.66665999999998
$ws.Range("N41").Value = -3153
$ws.Range("H70").Value = 8900
$ws.Range("I70").Value = 7700
$ws.Range("J70").Value = 9200
$ws.Range("K70").Value = 23100
$ws.Range("L70").Value = 27600
$ws.Range("M70").Value = -22830
$ws.Range("N70").Value = -28140
$ws.Range("H73").Value = 8900
$ws.Range("I73").Value = 7700
$ws.Range("J73").Value = 9200
$ws.Range("K73").Value = 23100
$ws.Range("L73").Value = 27600
$ws.Range("M73").Value = -22164
$ws.Range("N73").Value = -29472
$ws.Range("H86").Value = 142860820
$ws.Range("I86").Value = 333336220
$ws.Range("J86").Value = 4249.5
$ws.Range("K86").Value = 333336220
$ws.Range("L86").Value = 4249.5
$ws.Range("M86").Value = -333335097
$ws.Range("N86").Value = -6495.5
$ws.Range("H89").Value = 142860820
$ws.Range("I89").Value = 333336220
$ws.Range("J89").Value = 4249.5
$ws.Range("K89").Value = 1666681100
$ws.Range("L89").Value = 21247.5
$ws.Range("M89").Value = -1666675484
$ws.Range("N89").Value = -32479.5
$ws.Range("H132").Value = 18249.273
$ws.Range("I132").Value = 3458.6667
$ws.Range("K132").Value = 10376.0001
$ws.Range("M132").Value = -7846.000100000001
$ws.Range("H137").Value = 2995.9375
$ws.Range("I137").Value = 2802.7693
$ws.Range("K137").Value = 8408.3079
$ws.Range("M137").Value = -5858.3079

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12104.282
$ws.Range("I32").Value = 11385.397
$ws.Range("J32").Value = 22600
$ws.Range("K32").Value = 11385.397
$ws.Range("L32").Value = 22600
$ws.Range("M32").Value = -11098.397
$ws.Range("N32").Value = -23174
$ws.Range("H45").Value = 2628.5715
$ws.Range("J45").Value = 3614.2856
$ws.Range("L45").Value = 3614.2856
$ws.Range("N45").Value = -4368.2856
$ws.Range("H88").Value = 104741.4
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 104741.4
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 104741.4
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -105553.4
$ws.Range("H91").Value = 104741.4
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 104741.4
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 104741.4
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -107549.4
$ws.Range("H110").Value = 987.7646999999999
$ws.Range("I110").Value = 783.1429000000001
$ws.Range("K110").Value = 783.1429000000001
$ws.Range("M110").Value = 1261.8571

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4454.364
$ws.Range("I86").Value = 3862.2
$ws.Range("K86").Value = 3862.2
$ws.Range("M86").Value = -2739.2
$ws.Range("H89").Value = 4454.364
$ws.Range("I89").Value = 3862.2
$ws.Range("K89").Value = 19311
$ws.Range("M89").Value = -13695
$ws.Range("H107").Value = 2659.2307
$ws.Range("I107").Value = 2040.2
$ws.Range("K107").Value = 2040.2
$ws.Range("M107").Value = -120.2
$ws.Range("H132").Value = 92773.5
$ws.Range("J132").Value = 92773.5
$ws.Range("L132").Value = 92773.5
$ws.Range("N132").Value = -102893.5
$ws.Range("H134").Value = 5717.6665
$ws.Range("I134").Value = 5407.3335
$ws.Range("K134").Value = 16222.0005
$ws.Range("M134").Value = -13687.0005
$ws.Range("H135").Value = 104000
$ws.Range("J135").Value = 104000
$ws.Range("L135").Value = 104000
$ws.Range("N135").Value = -114140

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1048.9333
$ws.Range("I16").Value = 797.4
$ws.Range("K16").Value = 797.4
$ws.Range("M16").Value = -510.4
$ws.Range("H62").Value = 92000
$ws.Range("I62").Value = 50000
$ws.Range("K62").Value = 50000
$ws.Range("M62").Value = -49376
$ws.Range("H65").Value = 92000
$ws.Range("I65").Value = 50000
$ws.Range("K65").Value = 250000
$ws.Range("M65").Value = -246880
$ws.Range("H105").Value = 1896.75
$ws.Range("I105").Value = 1590
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 1590
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = 157
$ws.Range("N105").Value = -5493
$ws.Range("H107").Value = 724.4
$ws.Range("I107").Value = 629.1667
$ws.Range("J107").Value = 867.25
$ws.Range("K107").Value = 629.1667
$ws.Range("L107").Value = 867.25
$ws.Range("M107").Value = 1290.8333
$ws.Range("N107").Value = -4707.25
$ws.Range("H113").Value = 1048.9333
$ws.Range("I113").Value = 797.4
$ws.Range("K113").Value = 797.4
$ws.Range("M113").Value = 1372.6
$ws.Range("H140").Value = 113333.336
$ws.Range("J140").Value = 113333.336
$ws.Range("L140").Value = 113333.336
$ws.Range("N140").Value = -123693.336

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1322.875
$ws.Range("I34").Value = 352.22223
$ws.Range("J34").Value = 2570.8572
$ws.Range("K34").Value = 1056.66669
$ws.Range("L34").Value = 7712.571599999999
$ws.Range("M34").Value = -972.66669
$ws.Range("N34").Value = -7880.571599999999
$ws.Range("H39").Value = 4474.1665
$ws.Range("J39").Value = 4380
$ws.Range("L39").Value = 13140
$ws.Range("N39").Value = -13728
$ws.Range("H55").Value = 4396.1665
$ws.Range("J55").Value = 4522.727
$ws.Range("L55").Value = 13568.181
$ws.Range("N55").Value = -13922.181
$ws.Range("H122").Value = 970.2
$ws.Range("J122").Value = 916
$ws.Range("L122").Value = 8244
$ws.Range("N122").Value = -13144
$ws.Range("H133").Value = 9899.85
$ws.Range("I133").Value = 7593.5625
$ws.Range("J133").Value = 19125
$ws.Range("K133").Value = 22780.6875
$ws.Range("L133").Value = 57375
$ws.Range("M133").Value = -17720.6875
$ws.Range("N133").Value = -67495

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4749.5
$ws.Range("I113").Value = 4666
$ws.Range("K113").Value = 4666
$ws.Range("M113").Value = -2496
$ws.Range("H126").Value = 5612.4614
$ws.Range("I126").Value = 2193.6667
$ws.Range("K126").Value = 6581.000100000001
$ws.Range("M126").Value = -4111.000100000001
$ws.Range("H136").Value = 21657.732
$ws.Range("J136").Value = 21657.732
$ws.Range("L136").Value = 64973.196
$ws.Range("N136").Value = -70073.196

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 903.8823
$ws.Range("J22").Value = 896.8
$ws.Range("L22").Value = 896.8
$ws.Range("N22").Value = -1486.8
$ws.Range("H27").Value = 903.8823
$ws.Range("J27").Value = 896.8
$ws.Range("L27").Value = 896.8
$ws.Range("N27").Value = -1110.8
$ws.Range("H93").Value = 3599.6667
$ws.Range("I93").Value = 2899.5
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 2899.5
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -1651.5
$ws.Range("N93").Value = -7496
$ws.Range("H122").Value = 6130.3823
$ws.Range("I122").Value = 2978.5386
$ws.Range("J122").Value = 8081.524
$ws.Range("K122").Value = 8935.6158
$ws.Range("L122").Value = 24244.572
$ws.Range("M122").Value = -6485.6158
$ws.Range("N122").Value = -29144.572

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 53989.668
$ws.Range("J68").Value = 53989.668
$ws.Range("L68").Value = 53989.668
$ws.Range("N68").Value = -55611.668
$ws.Range("H71").Value = 53989.668
$ws.Range("J71").Value = 53989.668
$ws.Range("L71").Value = 161969.004
$ws.Range("N71").Value = -170081.004
$ws.Range("H100").Value = 1054.8823
$ws.Range("I100").Value = 534.38464
$ws.Range("K100").Value = 1068.76928
$ws.Range("M100").Value = -527.76928
$ws.Range("H113").Value = 855.63635
$ws.Range("I113").Value = 682
$ws.Range("K113").Value = 2046
$ws.Range("M113").Value = 124
$ws.Range("H126").Value = 1990.4
$ws.Range("I126").Value = 1990.4
$ws.Range("K126").Value = 5971.200000000001
$ws.Range("M126").Value = -3501.200000000001
